$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2750
$ws.Range("J40").Value = 2666.6667
$ws.Range("L40").Value = 2666.6667
$ws.Range("N40").Value = -3016.6667
$ws.Range("H86").Value = 3678.4285
$ws.Range("I86").Value = 3383
$ws.Range("J86").Value = 3900
$ws.Range("K86").Value = 3383
$ws.Range("L86").Value = 3900
$ws.Range("M86").Value = -2260
$ws.Range("N86").Value = -6146
$ws.Range("H89").Value = 3678.4285
$ws.Range("I89").Value = 3383
$ws.Range("J89").Value = 3900
$ws.Range("K89").Value = 16915
$ws.Range("L89").Value = 19500
$ws.Range("M89").Value = -11299
$ws.Range("N89").Value = -30732
$ws.Range("H111").Value = 2150
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3274.1333
$ws.Range("I88").Value = 3026.1667
$ws.Range("J88").Value = 3439.4443
$ws.Range("K88").Value = 3026.1667
$ws.Range("L88").Value = 3439.4443
$ws.Range("M88").Value = -2620.1667
$ws.Range("N88").Value = -4251.4443
$ws.Range("H91").Value = 3274.1333
$ws.Range("I91").Value = 3026.1667
$ws.Range("J91").Value = 3439.4443
$ws.Range("K91").Value = 3026.1667
$ws.Range("L91").Value = 3439.4443
$ws.Range("M91").Value = -1622.1667
$ws.Range("N91").Value = -6247.4443
$ws.Range("H94").Value = 59996
$ws.Range("J94").Value = 59996
$ws.Range("L94").Value = 59996
$ws.Range("N94").Value = -61798
$ws.Range("H102").Value = 2080.25
$ws.Range("I102").Value = 808.5
$ws.Range("J102").Value = 4199.8335
$ws.Range("K102").Value = 808.5
$ws.Range("L102").Value = 4199.8335
$ws.Range("M102").Value = 813.5
$ws.Range("N102").Value = -7443.8335
$ws.Range("H110").Value = 5993.5713
$ws.Range("I110").Value = 5993.5713
$ws.Range("K110").Value = 5993.5713
$ws.Range("M110").Value = -3948.5713

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2736.647
$ws.Range("I99").Value = 2672.8333
$ws.Range("K99").Value = 2672.8333
$ws.Range("M99").Value = -1174.8333
$ws.Range("H105").Value = 3601.3333
$ws.Range("I105").Value = 3121.6
$ws.Range("K105").Value = 3121.6
$ws.Range("M105").Value = -1374.6

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 625
$ws.Range("I8").Value = 600
$ws.Range("K8").Value = 600
$ws.Range("M8").Value = -460
$ws.Range("H16").Value = 961.75
$ws.Range("I16").Value = 915.6667
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 915.6667
$ws.Range("L16").Value = 1100
$ws.Range("M16").Value = -628.6667
$ws.Range("N16").Value = -1674
$ws.Range("H62").Value = 1966.3334
$ws.Range("I62").Value = 1849.5
$ws.Range("K62").Value = 1849.5
$ws.Range("M62").Value = -1225.5
$ws.Range("H65").Value = 1966.3334
$ws.Range("I65").Value = 1849.5
$ws.Range("K65").Value = 9247.5
$ws.Range("M65").Value = -6127.5
$ws.Range("H99").Value = 3782.8
$ws.Range("I99").Value = 3210.4375
$ws.Range("J99").Value = 4800.3335
$ws.Range("K99").Value = 3210.4375
$ws.Range("L99").Value = 4800.3335
$ws.Range("M99").Value = -1712.4375
$ws.Range("N99").Value = -7796.3335
$ws.Range("H103").Value = 6950
$ws.Range("I103").Value = 10000
$ws.Range("J103").Value = 3900
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 3900
$ws.Range("M103").Value = -8828
$ws.Range("N103").Value = -6244
$ws.Range("H105").Value = 1027.2609
$ws.Range("I105").Value = 648.7646999999999
$ws.Range("K105").Value = 648.7646999999999
$ws.Range("M105").Value = 1098.2353
$ws.Range("H113").Value = 961.75
$ws.Range("I113").Value = 915.6667
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 915.6667
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1254.3333
$ws.Range("N113").Value = -5440
$ws.Range("H126").Value = 3782.8
$ws.Range("I126").Value = 3210.4375
$ws.Range("J126").Value = 4800.3335
$ws.Range("K126").Value = 9631.3125
$ws.Range("L126").Value = 14401.0005
$ws.Range("M126").Value = -7161.3125
$ws.Range("N126").Value = -19341.0005
$ws.Range("H132").Value = 3799.4614
$ws.Range("I132").Value = 4103.8696
$ws.Range("J132").Value = 1465.6666
$ws.Range("K132").Value = 12311.6088
$ws.Range("L132").Value = 4396.9998
$ws.Range("M132").Value = -9781.6088
$ws.Range("N132").Value = -9456.9998

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1294.091
$ws.Range("I34").Value = 47.6
$ws.Range("J34").Value = 2332.8333
$ws.Range("K34").Value = 142.8
$ws.Range("L34").Value = 6998.499899999999
$ws.Range("M34").Value = -58.80000000000001
$ws.Range("N34").Value = -7166.499899999999
$ws.Range("H104").Value = 1000
$ws.Range("J104").Value = 1000
$ws.Range("L104").Value = 3000
$ws.Range("N104").Value = -8242
$ws.Range("H122").Value = 1975
$ws.Range("J122").Value = 1974.6666
$ws.Range("L122").Value = 17771.9994
$ws.Range("N122").Value = -22671.9994

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H70").Value = 7697.857
$ws.Range("I70").Value = 7697.857
$ws.Range("K70").Value = 7697.857
$ws.Range("M70").Value = -7427.857
$ws.Range("H73").Value = 7697.857
$ws.Range("I73").Value = 7697.857
$ws.Range("K73").Value = 7697.857
$ws.Range("M73").Value = -6761.857
$ws.Range("H132").Value = 2964.6667
$ws.Range("I132").Value = 2186
$ws.Range("J132").Value = 4989.2
$ws.Range("K132").Value = 6558
$ws.Range("L132").Value = 14967.6
$ws.Range("M132").Value = -4028
$ws.Range("N132").Value = -20027.6

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 1336.6666
$ws.Range("I10").Value = 282.5
$ws.Range("J10").Value = 3445
$ws.Range("K10").Value = 282.5
$ws.Range("L10").Value = 3445
$ws.Range("M10").Value = -142.5
$ws.Range("N10").Value = -3725
$ws.Range("H17").Value = 15999.5
$ws.Range("I17").Value = 9999
$ws.Range("J17").Value = 22000
$ws.Range("K17").Value = 9999
$ws.Range("L17").Value = 22000
$ws.Range("M17").Value = -9829
$ws.Range("N17").Value = -22340
$ws.Range("H22").Value = 1342.1428
$ws.Range("J22").Value = 1166.6666
$ws.Range("L22").Value = 1166.6666
$ws.Range("N22").Value = -1756.6666
$ws.Range("H27").Value = 1342.1428
$ws.Range("J27").Value = 1166.6666
$ws.Range("L27").Value = 1166.6666
$ws.Range("N27").Value = -1380.6666
$ws.Range("H55").Value = 718.4231
$ws.Range("I55").Value = 163.85715
$ws.Range("J55").Value = 922.7368
$ws.Range("K55").Value = 163.85715
$ws.Range("L55").Value = 922.7368
$ws.Range("M55").Value = 9.14285000000001
$ws.Range("N55").Value = -1268.7368
$ws.Range("H61").Value = 2661.25
$ws.Range("I61").Value = 2661.25
$ws.Range("K61").Value = 2661.25
$ws.Range("M61").Value = -2459.25
$ws.Range("H82").Value = 2761.9375
$ws.Range("J82").Value = 3308.0833
$ws.Range("L82").Value = 3308.0833
$ws.Range("N82").Value = -4030.0833
$ws.Range("H85").Value = 2761.9375
$ws.Range("J85").Value = 3308.0833
$ws.Range("L85").Value = 3308.0833
$ws.Range("N85").Value = -5804.0833
$ws.Range("H113").Value = 2661.25
$ws.Range("I113").Value = 2661.25
$ws.Range("K113").Value = 2661.25
$ws.Range("M113").Value = -491.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 597.5
$ws.Range("J6").Value = 597.5
$ws.Range("L6").Value = 597.5
$ws.Range("N6").Value = -827.5
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H132").Value = 1267.6923
$ws.Range("I132").Value = 1238.4
$ws.Range("K132").Value = 3715.2
$ws.Range("M132").Value = -1185.2
